$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: permutation of 0..19 representing the population "sequence to use"
$ws.Range("B1").Value = 3
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 1
$ws.Range("B4").Value = 13
$ws.Range("B5").Value = 16
$ws.Range("B6").Value = 6
$ws.Range("B7").Value = 15
$ws.Range("B8").Value = 5
$ws.Range("B9").Value = 19
$ws.Range("B10").Value = 8
$ws.Range("B11").Value = 2
$ws.Range("B12").Value = 9
$ws.Range("B13").Value = 11
$ws.Range("B14").Value = 4
$ws.Range("B15").Value = 17
$ws.Range("B16").Value = 7
$ws.Range("B17").Value = 12
$ws.Range("B18").Value = 14
$ws.Range("B19").Value = 18
$ws.Range("B20").Value = 10

# Column D: fitness / penalty improvement values
$ws.Range("D1").Value = 109.9919950107364
$ws.Range("D2").Value = 74.44993474229956

# Last generation fit
$ws.Range("B21").Value = 0.7906977336495797
